# Salt_Lake_Data.xlsx edit:
# Two more "Unnamed: 0...." index columns were re-introduced into the source
# data (as if an extra reset_index() had been applied twice), which inserts
# two new numeric columns right before the "Date" column. Everything from
# the old "Date" column onward shifts two columns to the right
# (O..S -> Q..U), and the two freshly inserted columns (new O, P) are filled
# with the same row-index value that already populates columns A..N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank columns at O:P. Everything at O:S shifts to Q:U,
#    carrying its data/headers/formatting along for the ride.
$ws.Range("O:P").Insert()

# 2) New header text for the two freshly-inserted (now-empty) columns,
#    matching the bold/border/centered header formatting already used by
#    the rest of the header row (column N is a representative sample).
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("O1").Value2 = "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1.1.1"
$ws.Range("P1").Value2 = "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1.1.1.1"

# 3) Fill the new O/P columns with the row-index value already present in
#    column A (and every other index column, A..N) for each data row.
for ($r = 2; $r -le 83; $r++) {
    $idxVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 15).Value2 = $idxVal
    $ws.Cells.Item($r, 16).Value2 = $idxVal
}

# 4) Row 84 is sparse in the source data: only column A (and O..S) were
#    populated, columns B..N were never filled in. The same two new
#    index columns there land on the first two still-empty cells in the
#    row, which are B84/C84 rather than O84/P84.
$idxVal84 = $ws.Cells.Item(84, 1).Value2
$ws.Cells.Item(84, 2).Value2 = $idxVal84
$ws.Cells.Item(84, 3).Value2 = $idxVal84
